# Update the "想去人数" (want-to-go count) figures in column F for the
# sheets that contain the scraped convention data: "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1577
    "F4"  = 1033
    "F7"  = 2687
    "F9"  = 1713
    "F11" = 73
    "F12" = 579
    "F14" = 13
    "F15" = 89
    "F16" = 76
    "F17" = 81
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
